# Apply "ER atualizado e tabelas mapeadas" edits to the "user" mapping sheet.
#
# Summary of the change:
#  - Rows 20 and 31 (cidade/city_id and function_id/function_id) are marked
#    as already-reviewed by giving them a white ("Background 1") fill and
#    clearing the helper column C note.
#  - Two new rows are appended (34, 35) documenting new FK columns that were
#    added to the ER model: "XXX" -> profile_id and "XXX" -> platoon_id,
#    highlighted in yellow to call out that they are new/unmapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark rows 20 and 31 as reviewed (white fill, column C note removed) ---
foreach ($r in 20, 31) {
    $rowRange = $ws.Range("A" + $r + ":B" + $r)
    $rowRange.Interior.ThemeColor = 2   # xlThemeColorLight1 -> theme="0" (white)
    $ws.Cells.Item($r, 3).ClearContents()
}

# --- Append the two new mapped columns discovered in the updated ER ---
$ws.Range("A34").Value = "XXX"
$ws.Range("B34").Value = "profile_id"

$ws.Range("A35").Value = "XXX"
$ws.Range("B35").Value = "platoon_id"

$newRows = $ws.Range("A34:B35")
$newRows.Interior.Color = 65535   # RGB(255,255,0) -> yellow highlight

Write-Output "done"
